# Trade #90 closed at 2026-02-17 09:08:42 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up numbers to reflect the new
# closed trade, and appends the new trade row to both the "All Trades" and
# "MarketMaking" detail sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.69
$summary.Range("B4").Value = -0.3
$summary.Range("B5").Value = -0.07000000000000001
$summary.Range("B6").Value = 90
$summary.Range("B8").Value = 36
$summary.Range("B9").Value = 42.22

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.69
$status.Range("D4").Value = 90
$status.Range("E4").Value = -0.3
$status.Range("F4").Value = -0.31
$status.Range("G4").Value = 42.22

# ---------------------------------------------------------------------
# New trade row data (trade #90)
# ---------------------------------------------------------------------
$newRow = 91

function Add-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value = 90

    # The date text looks numeric-ish, so force it to stay plain text
    # (matching the other rows, which store dates as literal strings)
    # instead of letting Excel auto-convert it into a date serial.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "09:08:36"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.87
    $ws.Cells.Item($row, 7).Value = 0.82
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -5.7471
    $ws.Cells.Item($row, 10).Value = -0.05
    $ws.Cells.Item($row, 11).Value = 99.69
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.1
}

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades $newRow

# ---------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking $newRow
